$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifts old row 5 "sCs" target down to row 6)
$ws.Rows.Item(5).Insert()

# Update row 2 (ECs target) - ligand stats and receptor stats
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.199962333333334
$ws.Range("H2").Value = 12.599887
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.314108
$ws.Range("N2").Value = 3.942324
$ws.Range("O2").Value = 0.05768654525237047
$ws.Range("P2").Value = 0.05768654525237048
$ws.Range("Q2").Value = 5.519204101932
$ws.Range("R2").Value = 49.672836917388
$ws.Range("S2").Value = 0.05768654525237047
$ws.Range("T2").Value = 0.05768654525237048

# Update row 3 (FAPs target)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.199962333333334
$ws.Range("H3").Value = 12.599887
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.912216333333333
$ws.Range("N3").Value = 14.736649
$ws.Range("O3").Value = 0.2156358456095441
$ws.Range("P3").Value = 0.2156358456095441
$ws.Range("Q3").Value = 20.63112357318478
$ws.Range("R3").Value = 185.680112158663
$ws.Range("S3").Value = 0.2156358456095441
$ws.Range("T3").Value = 0.2156358456095441

# Update row 4 - Target cluster changes from Neutro to M1 (new cluster inserted before Neutro)
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.199962333333334
$ws.Range("H4").Value = 12.599887
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0405
$ws.Range("N4").Value = 0.1215
$ws.Range("O4").Value = 0.001777863830614382
$ws.Range("P4").Value = 0.001777863830614382
$ws.Range("Q4").Value = 0.1700984745
$ws.Range("R4").Value = 1.5308862705
$ws.Range("S4").Value = 0.001777863830614382
$ws.Range("T4").Value = 0.001777863830614382

# Fill in new row 5 (Neutro target)
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Edn3"
$ws.Range("C5").Value = "Ednra"
$ws.Range("D5").Value = "Neutro"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.199962333333334
$ws.Range("H5").Value = 12.599887
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4069286666666667
$ws.Range("N5").Value = 1.220786
$ws.Range("O5").Value = 0.01786330266930378
$ws.Range("P5").Value = 0.01786330266930378
$ws.Range("Q5").Value = 1.709085072353556
$ws.Range("R5").Value = 15.381765651182
$ws.Range("S5").Value = 0.01786330266930378
$ws.Range("T5").Value = 0.01786330266930378

# Update row 6 (was row 5, sCs target) with new values
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.199962333333334
$ws.Range("H6").Value = 12.599887
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 16.10639433333333
$ws.Range("N6").Value = 48.319183
$ws.Range("O6").Value = 0.7070364426381672
$ws.Range("P6").Value = 0.7070364426381672
$ws.Range("Q6").Value = 67.64624952581343
$ws.Range("R6").Value = 608.816245732321
$ws.Range("S6").Value = 0.7070364426381672
$ws.Range("T6").Value = 0.7070364426381672
